$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update the CodeSystem URL (pythia -> cicada)
$ws.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/EvalReason"

# 2. Update the generation Date
$ws.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# 3. Insert a new "Jurisdiction" metadata row after "Contact" (row 10), before "Description"
$ws.Rows.Item(11).Insert()

# Match the formatting of the surrounding metadata rows for the new row
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
